# Apply edits to the Raw_Annotations sheet:
# - Set E2:E31 to 300 (formulas in J/K recompute automatically)
# - Update sheet view: normalize zoom to 100%,
#   change selection to E2:E31 with active cell E2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw_Annotations")
$ws.Activate()

# Update the E column (travel time / headway, etc.) for rows 2-31 to 300
$ws.Range("E2:E31").Value = 300

# Reset the zoom level of the sheet view to 100%
$excel.ActiveWindow.Zoom = 100

# Update the selection / view to match the new focus on column E
$ws.Range("E2:E31").Select()

# Recalculate the workbook so dependent formulas (J, K columns) update
$excel.CalculateFullRebuild()
